$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 314738
$ws.Cells.Item(2, 4).Value = 401280715
$ws.Cells.Item(4, 3).Value = 313
$ws.Cells.Item(4, 4).Value = 447707
$ws.Cells.Item(6, 3).Value = 8
$ws.Cells.Item(6, 4).Value = 11676
$ws.Cells.Item(8, 3).Value = 845
$ws.Cells.Item(8, 4).Value = 1245107
$ws.Cells.Item(10, 3).Value = 115723
$ws.Cells.Item(10, 4).Value = 169581339
$ws.Cells.Item(12, 3).Value = 58345
$ws.Cells.Item(12, 4).Value = 84213662
$ws.Cells.Item(14, 3).Value = 48
$ws.Cells.Item(14, 4).Value = 66043
$ws.Cells.Item(16, 3).Value = 3939
$ws.Cells.Item(16, 4).Value = 5592761
$ws.Cells.Item(20, 3).Value = 6414
$ws.Cells.Item(20, 4).Value = 8951081
$ws.Cells.Item(22, 3).Value = 76141
$ws.Cells.Item(22, 4).Value = 95044030
$ws.Cells.Item(28, 3).Value = 32127
$ws.Cells.Item(28, 4).Value = 47036600
$ws.Cells.Item(30, 3).Value = 11315
$ws.Cells.Item(30, 4).Value = 16277294
$ws.Cells.Item(33, 3).Value = 1545
$ws.Cells.Item(33, 4).Value = 2168807
$ws.Cells.Item(35, 3).Value = 1763
$ws.Cells.Item(35, 4).Value = 2487523
$ws.Cells.Item(36, 3).Value = 95685
$ws.Cells.Item(36, 4).Value = 120561110
$ws.Cells.Item(42, 3).Value = 897
$ws.Cells.Item(42, 4).Value = 1320185
$ws.Cells.Item(44, 3).Value = 43970
$ws.Cells.Item(44, 4).Value = 64451538
$ws.Cells.Item(45, 3).Value = 25
$ws.Cells.Item(45, 4).Value = 37450
$ws.Cells.Item(46, 3).Value = 8997
$ws.Cells.Item(46, 4).Value = 12914224
$ws.Cells.Item(48, 3).Value = 1383
$ws.Cells.Item(48, 4).Value = 1920900
$ws.Cells.Item(51, 3).Value = 2233
$ws.Cells.Item(51, 4).Value = 3114064
$ws.Cells.Item(52, 3).Value = 67884
$ws.Cells.Item(52, 4).Value = 85214566
$ws.Cells.Item(58, 3).Value = 27826
$ws.Cells.Item(58, 4).Value = 40811711
$ws.Cells.Item(61, 3).Value = 10917
$ws.Cells.Item(61, 4).Value = 15784890
$ws.Cells.Item(67, 3).Value = 1420
$ws.Cells.Item(67, 4).Value = 1987268
$ws.Cells.Item(69, 3).Value = 20144
$ws.Cells.Item(69, 4).Value = 26385968
$ws.Cells.Item(73, 3).Value = 7467
$ws.Cells.Item(73, 4).Value = 10932361
$ws.Cells.Item(75, 3).Value = 5023
$ws.Cells.Item(75, 4).Value = 7292706
$ws.Cells.Item(76, 3).Value = 478
$ws.Cells.Item(76, 4).Value = 675739
$ws.Cells.Item(78, 3).Value = 138296
$ws.Cells.Item(78, 4).Value = 172530287
$ws.Cells.Item(84, 3).Value = 62830
$ws.Cells.Item(84, 4).Value = 92094575
$ws.Cells.Item(87, 3).Value = 29235
$ws.Cells.Item(87, 4).Value = 42298498
$ws.Cells.Item(89, 3).Value = 2689
$ws.Cells.Item(89, 4).Value = 3872520
$ws.Cells.Item(90, 3).Value = 2723
$ws.Cells.Item(90, 4).Value = 3847605
$ws.Cells.Item(91, 3).Value = 31684
$ws.Cells.Item(91, 4).Value = 42935818
$ws.Cells.Item(95, 3).Value = 7732
$ws.Cells.Item(95, 4).Value = 11369567
$ws.Cells.Item(97, 3).Value = 7027
$ws.Cells.Item(97, 4).Value = 10186733
$ws.Cells.Item(99, 3).Value = 511
$ws.Cells.Item(99, 4).Value = 728005
$ws.Cells.Item(100, 3).Value = 477
$ws.Cells.Item(100, 4).Value = 688443
$ws.Cells.Item(101, 3).Value = 8565
$ws.Cells.Item(101, 4).Value = 11884856
$ws.Cells.Item(103, 3).Value = 2155
$ws.Cells.Item(103, 4).Value = 3174470
$ws.Cells.Item(105, 3).Value = 2905
$ws.Cells.Item(105, 4).Value = 4241402
$ws.Cells.Item(107, 3).Value = 118
$ws.Cells.Item(107, 4).Value = 171120
$ws.Cells.Item(108, 3).Value = 159
$ws.Cells.Item(108, 4).Value = 224586
$ws.Cells.Item(109, 3).Value = 138837
$ws.Cells.Item(109, 4).Value = 171727843
$ws.Cells.Item(115, 3).Value = 52110
$ws.Cells.Item(115, 4).Value = 76395826
$ws.Cells.Item(116, 3).Value = 85
$ws.Cells.Item(116, 4).Value = 125959
$ws.Cells.Item(117, 3).Value = 26491
$ws.Cells.Item(117, 4).Value = 38379200
$ws.Cells.Item(118, 3).Value = 1294
$ws.Cells.Item(118, 4).Value = 1770551
$ws.Cells.Item(121, 3).Value = 2175
$ws.Cells.Item(121, 4).Value = 3054218
$ws.Cells.Item(123, 3).Value = 490570
$ws.Cells.Item(123, 4).Value = 647008988
$ws.Cells.Item(125, 3).Value = 208
$ws.Cells.Item(125, 4).Value = 306736
$ws.Cells.Item(128, 3).Value = 1355
$ws.Cells.Item(128, 4).Value = 2008811
$ws.Cells.Item(130, 3).Value = 203791
$ws.Cells.Item(130, 4).Value = 299590510
$ws.Cells.Item(131, 3).Value = 386
$ws.Cells.Item(131, 4).Value = 575790
$ws.Cells.Item(133, 3).Value = 175987
$ws.Cells.Item(133, 4).Value = 255802165
$ws.Cells.Item(136, 3).Value = 2793
$ws.Cells.Item(136, 4).Value = 3925444
$ws.Cells.Item(138, 3).Value = 6121
$ws.Cells.Item(138, 4).Value = 8647603
$ws.Cells.Item(141, 3).Value = 43553
$ws.Cells.Item(141, 4).Value = 58160226
$ws.Cells.Item(147, 3).Value = 13836
$ws.Cells.Item(147, 4).Value = 20294557
$ws.Cells.Item(148, 3).Value = 3685
$ws.Cells.Item(148, 4).Value = 5315965
$ws.Cells.Item(151, 3).Value = 389
$ws.Cells.Item(151, 4).Value = 559431
$ws.Cells.Item(154, 3).Value = 17125
$ws.Cells.Item(154, 4).Value = 22633593
$ws.Cells.Item(158, 3).Value = 7001
$ws.Cells.Item(158, 4).Value = 10185416
$ws.Cells.Item(160, 3).Value = 4874
$ws.Cells.Item(160, 4).Value = 7015456
$ws.Cells.Item(163, 3).Value = 258
$ws.Cells.Item(163, 4).Value = 369783
$ws.Cells.Item(165, 3).Value = 14796
$ws.Cells.Item(165, 4).Value = 21463720
$ws.Cells.Item(166, 3).Value = 1716
$ws.Cells.Item(166, 4).Value = 2552630
$ws.Cells.Item(167, 3).Value = 232
$ws.Cells.Item(167, 4).Value = 342802
$ws.Cells.Item(169, 3).Value = 48
$ws.Cells.Item(169, 4).Value = 71690
$ws.Cells.Item(171, 3).Value = 86300
$ws.Cells.Item(171, 4).Value = 107985744
$ws.Cells.Item(176, 3).Value = 638
$ws.Cells.Item(176, 4).Value = 940348
$ws.Cells.Item(178, 3).Value = 33473
$ws.Cells.Item(178, 4).Value = 49091394
$ws.Cells.Item(180, 3).Value = 12789
$ws.Cells.Item(180, 4).Value = 18478490
$ws.Cells.Item(182, 3).Value = 1231
$ws.Cells.Item(182, 4).Value = 1722396
$ws.Cells.Item(184, 3).Value = 1594
$ws.Cells.Item(184, 4).Value = 2244193
$ws.Cells.Item(186, 3).Value = 234202
$ws.Cells.Item(186, 4).Value = 291214392
$ws.Cells.Item(194, 3).Value = 85668
$ws.Cells.Item(194, 4).Value = 125590329
$ws.Cells.Item(197, 3).Value = 32501
$ws.Cells.Item(197, 4).Value = 46776170
$ws.Cells.Item(200, 3).Value = 5022
$ws.Cells.Item(200, 4).Value = 7156762
$ws.Cells.Item(203, 3).Value = 4691
$ws.Cells.Item(203, 4).Value = 6487816
$ws.Cells.Item(206, 3).Value = 258820
$ws.Cells.Item(206, 4).Value = 320397886
$ws.Cells.Item(208, 3).Value = 247
$ws.Cells.Item(208, 4).Value = 353087
$ws.Cells.Item(215, 3).Value = 93974
$ws.Cells.Item(215, 4).Value = 137494162
$ws.Cells.Item(218, 3).Value = 50574
$ws.Cells.Item(218, 4).Value = 73092217
$ws.Cells.Item(221, 3).Value = 4597
$ws.Cells.Item(221, 4).Value = 6452415
$ws.Cells.Item(224, 3).Value = 5525
$ws.Cells.Item(224, 4).Value = 7639436
$ws.Cells.Item(227, 3).Value = 104441
$ws.Cells.Item(227, 4).Value = 130749202
$ws.Cells.Item(232, 3).Value = 561
$ws.Cells.Item(232, 4).Value = 819439
$ws.Cells.Item(234, 3).Value = 48959
$ws.Cells.Item(234, 4).Value = 71731859
$ws.Cells.Item(236, 3).Value = 12173
$ws.Cells.Item(236, 4).Value = 17500151
$ws.Cells.Item(238, 3).Value = 1875
$ws.Cells.Item(238, 4).Value = 2687338
$ws.Cells.Item(240, 3).Value = 2421
$ws.Cells.Item(240, 4).Value = 3382219
$ws.Cells.Item(241, 3).Value = 252635
$ws.Cells.Item(241, 4).Value = 319104043
$ws.Cells.Item(244, 3).Value = 15
$ws.Cells.Item(244, 4).Value = 22500
$ws.Cells.Item(247, 3).Value = 819
$ws.Cells.Item(247, 4).Value = 1203050
$ws.Cells.Item(249, 3).Value = 94559
$ws.Cells.Item(249, 4).Value = 138570298
$ws.Cells.Item(252, 3).Value = 63786
$ws.Cells.Item(252, 4).Value = 92438746
$ws.Cells.Item(254, 3).Value = 2369
$ws.Cells.Item(254, 4).Value = 3343918
$ws.Cells.Item(257, 3).Value = 4453
$ws.Cells.Item(257, 4).Value = 6250224
